$d = $word.ActiveDocument

# The 2nd paragraph currently holds the old "Planteamiento del problema: ..."
# run (Segoe UI / shaded). The target revision:
#   1. Inserts a brand-new (plain, no paragraph formatting) paragraph right
#      before it with the new wording, styled with the "s1ppyq" character
#      style + black color.
#   2. Empties that old paragraph's run, leaving only its (unchanged) empty
#      paragraph mark / pPr behind.
#
# InsertXML on the old paragraph's Range does both at once: it replaces the
# range's contents (the old run text) with our fragment, which itself is a
# lone <w:p> holding the new run - Word inserts it as a new paragraph ahead
# of the (now emptied) original paragraph mark.
$oldPara = $d.Paragraphs.Item(2)
$oldRange = $oldPara.Range
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">El restaurante La pescadería ha identificado que el proceso de toma de pedidos y entrega de alimentos a los clientes puede resultar lento y poco eficiente durante horas pico. Esto se debe a la gran cantidad de comensales que acuden al establecimiento y al proceso manual de toma de pedidos, que a menudo conduce a errores y retrasos. ​</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$oldRange.InsertXML($xmlFrag) | Out-Null

# Apply the character style (rStyle) + explicit black color to the run of the
# newly inserted paragraph (now paragraph #2). The sub-range deliberately
# excludes the trailing paragraph mark so Word records "s1ppyq" as a run-level
# character style (w:rStyle) rather than promoting it to a paragraph style
# (w:pStyle).
$newPara = $d.Paragraphs.Item(2)
$newRange = $newPara.Range
$subRange = $d.Range($newRange.Start, $newRange.End - 1)
$subRange.Style = "s1ppyq"
$subRange.Font.Color = 0
